$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 203.7816646666667
$ws.Range("H2").Value = 611.344994
$ws.Range("I2").Value = 0.6667327591988204
$ws.Range("J2").Value = 0.6667327591988205
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 2.08532
$ws.Range("N2").Value = 6.25596
$ws.Range("O2").Value = 0.01753772176136817
$ws.Range("P2").Value = 0.01753772176136816
$ws.Range("Q2").Value = 424.9499809626933
$ws.Range("R2").Value = 3824.54982866424
$ws.Range("S2").Value = 0.01169297362001819
$ws.Range("T2").Value = 0.01169297362001819

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 203.7816646666667
$ws.Range("H3").Value = 611.344994
$ws.Range("I3").Value = 0.6667327591988204
$ws.Range("J3").Value = 0.6667327591988205
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 101.898173
$ws.Range("N3").Value = 305.694519
$ws.Range("O3").Value = 0.8569724579756384
$ws.Range("P3").Value = 0.8569724579756383
$ws.Range("Q3").Value = 20764.97932043199
$ws.Range("R3").Value = 186884.8138838879
$ws.Range("S3").Value = 0.5713716114634925
$ws.Range("T3").Value = 0.5713716114634925

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 203.7816646666667
$ws.Range("H4").Value = 611.344994
$ws.Range("I4").Value = 0.6667327591988204
$ws.Range("J4").Value = 0.6667327591988205
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 14.921347
$ws.Range("N4").Value = 44.76404100000001
$ws.Range("O4").Value = 0.1254898202629935
$ws.Range("P4").Value = 0.1254898202629935
$ws.Range("Q4").Value = 3040.696930728973
$ws.Range("R4").Value = 27366.27237656076
$ws.Range("S4").Value = 0.0836681741153097
$ws.Range("T4").Value = 0.0836681741153097

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 63.14058933333333
$ws.Range("H5").Value = 189.421768
$ws.Range("I5").Value = 0.2065833519051582
$ws.Range("J5").Value = 0.2065833519051582
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.08532
$ws.Range("N5").Value = 6.25596
$ws.Range("O5").Value = 0.01753772176136817
$ws.Range("P5").Value = 0.01753772176136816
$ws.Range("Q5").Value = 131.6683337485867
$ws.Range("R5").Value = 1185.01500373728
$ws.Range("S5").Value = 0.00362300134624347
$ws.Range("T5").Value = 0.00362300134624347

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 63.14058933333333
$ws.Range("H6").Value = 189.421768
$ws.Range("I6").Value = 0.2065833519051582
$ws.Range("J6").Value = 0.2065833519051582
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 101.898173
$ws.Range("N6").Value = 305.694519
$ws.Range("O6").Value = 0.8569724579756384
$ws.Range("P6").Value = 0.8569724579756383
$ws.Range("Q6").Value = 6433.910695209955
$ws.Range("R6").Value = 57905.19625688959
$ws.Range("S6").Value = 0.1770362428590097
$ws.Range("T6").Value = 0.1770362428590097

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 63.14058933333333
$ws.Range("H7").Value = 189.421768
$ws.Range("I7").Value = 0.2065833519051582
$ws.Range("J7").Value = 0.2065833519051582
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 14.921347
$ws.Range("N7").Value = 44.76404100000001
$ws.Range("O7").Value = 0.1254898202629935
$ws.Range("P7").Value = 0.1254898202629935
$ws.Range("Q7").Value = 942.1426432271654
$ws.Range("R7").Value = 8479.283789044488
$ws.Range("S7").Value = 0.02592410769990504
$ws.Range("T7").Value = 0.02592410769990504

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 38.719942
$ws.Range("H8").Value = 116.159826
$ws.Range("I8").Value = 0.1266838888960214
$ws.Range("J8").Value = 0.1266838888960214
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 2.08532
$ws.Range("N8").Value = 6.25596
$ws.Range("O8").Value = 0.01753772176136817
$ws.Range("P8").Value = 0.01753772176136816
$ws.Range("Q8").Value = 80.74346945143999
$ws.Range("R8").Value = 726.6912250629599
$ws.Range("S8").Value = 0.002221746795106501
$ws.Range("T8").Value = 0.002221746795106501

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 38.719942
$ws.Range("H9").Value = 116.159826
$ws.Range("I9").Value = 0.1266838888960214
$ws.Range("J9").Value = 0.1266838888960214
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 101.898173
$ws.Range("N9").Value = 305.694519
$ws.Range("O9").Value = 0.8569724579756384
$ws.Range("P9").Value = 0.8569724579756383
$ws.Range("Q9").Value = 3945.491348465966
$ws.Range("R9").Value = 35509.42213619369
$ws.Range("S9").Value = 0.1085646036531361
$ws.Range("T9").Value = 0.1085646036531361

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 38.719942
$ws.Range("H10").Value = 116.159826
$ws.Range("I10").Value = 0.1266838888960214
$ws.Range("J10").Value = 0.1266838888960214
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 14.921347
$ws.Range("N10").Value = 44.76404100000001
$ws.Range("O10").Value = 0.1254898202629935
$ws.Range("P10").Value = 0.1254898202629935
$ws.Range("Q10").Value = 577.753690401874
$ws.Range("R10").Value = 5199.783213616866
$ws.Range("S10").Value = 0.01589753844777876
$ws.Range("T10").Value = 0.01589753844777876

